# Refresh the cryptos list (prices / 1h volume %) with the latest scrape.
# NumberFormat="@" is set before assigning any Price (column D) value that
# looks numeric so Excel stores it as text (matching the source feed, which
# renders things like "30.408.29" or "0.00001160" as literal strings) rather
# than silently coercing it to a Double and mangling the digits / formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.408.29"
$ws.Range("E2").Value = "  +2.23%  "
$ws.Range("D3").Value = "2.095.68"
$ws.Range("E3").Value = "  -0.14%  "
$ws.Range("E4").Value = "  -0.68%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "343.44"
$ws.Range("E6").Value = "  -0.62%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5275"
$ws.Range("E7").Value = "  +2.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4428"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.66"
$ws.Range("E9").Value = "  +3.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09343"
$ws.Range("E10").Value = "  +0.43%  "
$ws.Range("E11").Value = "  +0.36%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.78"
$ws.Range("E12").Value = "  -0.76%  "

# Rows 13-15 rotate: Chainlink moves up to 13, Polkadot up to 14, and
# WrappedEther (previously 13) drops to 15 - each with refreshed price/volume.
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.576"
$ws.Range("E13").Value = "  +3.73%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.919"
$ws.Range("E14").Value = "  +2.09%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "2.069.62"
$ws.Range("E15").Value = "  -1.79%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "101.55"
$ws.Range("E16").Value = "  +1.87%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001160"
$ws.Range("E17").Value = "  +0.44%  "
$ws.Range("E18").Value = "  -0.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "21.21"
$ws.Range("E19").Value = "  +1.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06693"
$ws.Range("E20").Value = "  +0.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.337"
$ws.Range("E21").Value = "  +2.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.003"
$ws.Range("E22").Value = "  -0.50%  "
$ws.Range("D23").Value = "30.435.45"
$ws.Range("E23").Value = "  +2.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.56"
$ws.Range("E24").Value = "  +0.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.319"
$ws.Range("E25").Value = "  +0.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "21.87"
$ws.Range("E26").Value = "  -0.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "163.13"
$ws.Range("E27").Value = "  +1.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.804"
$ws.Range("E28").Value = "  +8.28%  "
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.69"
$ws.Range("E30").Value = "  +0.40%  "
$ws.Range("E31").Value = "  -0.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1052"
$ws.Range("E32").Value = "  +0.21%  "
$ws.Range("E33").Value = "  +0.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.272"
$ws.Range("E34").Value = "  +1.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.868"
$ws.Range("E35").Value = "  -1.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.16"
$ws.Range("E36").Value = "  -0.88%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02633"
$ws.Range("E37").Value = "  +2.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06802"
$ws.Range("E38").Value = "  +1.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.7008"
$ws.Range("E39").Value = "  +1.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.61"
$ws.Range("E41").Value = "  +2.42%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2221"
$ws.Range("E42").Value = "  -0.61%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6873"
$ws.Range("E43").Value = "  +0.78%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.36"
$ws.Range("E44").Value = "  -0.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.342"
$ws.Range("E45").Value = "  +0.73%  "
$ws.Range("E46").Value = "  -0.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.390"
$ws.Range("E47").Value = "  +19.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.641"
$ws.Range("E48").Value = "  +0.41%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.240"
$ws.Range("E49").Value = "  +10.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00000000347"
$ws.Range("E50").Value = "  -2.51%  "
$ws.Range("E51").Value = "  -0.19%  "
